$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; E = 4; F = 20251215 }
    @{ Row = 3; E = 4; F = 20251215 }
    @{ Row = 4; E = 4; F = 20251215 }
    @{ Row = 5; E = 10; F = 20251225 }
    @{ Row = 6; E = 4; F = 20251215 }
    @{ Row = 7; E = 10; F = 20251225 }
    @{ Row = 8; E = 4; F = 20251215 }
    @{ Row = 9; E = 10; F = 20251225 }
    @{ Row = 10; E = 4; F = 20251222 }
    @{ Row = 11; E = 4; F = 20251215 }
    @{ Row = 12; E = 10; F = 20251225 }
    @{ Row = 13; E = 4; F = 20251215 }
    @{ Row = 14; E = 4; F = 20251215 }
    @{ Row = 15; E = 4; F = 20251215 }
    @{ Row = 16; E = 4; F = 20251219 }
    @{ Row = 17; E = 10; F = 20251225 }
    @{ Row = 18; E = 3; F = 20251218 }
    @{ Row = 19; E = 3; F = 20251218 }
    @{ Row = 20; E = 3; F = 20251218 }
    @{ Row = 21; E = 3; F = 20251218 }
    @{ Row = 22; E = 10; F = 20251225 }
    @{ Row = 23; E = 10; F = 20251225 }
    @{ Row = 24; E = 10; F = 20251225 }
    @{ Row = 25; E = 10; F = 20251225 }
    @{ Row = 26; E = 10; F = 20251225 }
    @{ Row = 27; E = 5; F = 20251223 }
    @{ Row = 28; E = 3; F = 20251218 }
    @{ Row = 29; E = 3; F = 20251218 }
    @{ Row = 30; E = 3; F = 20251218 }
    @{ Row = 31; E = 3; F = 20251218 }
    @{ Row = 32; E = 3; F = 20251218 }
    @{ Row = 33; E = 3; F = 20251218 }
    @{ Row = 34; E = 3; F = 20251218 }
    @{ Row = 35; E = 3; F = 20251218 }
    @{ Row = 37; E = 3; F = 20251218 }
    @{ Row = 38; E = 3; F = 20251218 }
    @{ Row = 39; E = 3; F = 20251218 }
    @{ Row = 40; E = 4; F = 20251222 }
    @{ Row = 41; E = 4; F = 20251222 }
    @{ Row = 42; E = 3; F = 20251218 }
    @{ Row = 43; E = 10; F = 20251225 }
    @{ Row = 44; E = 4; F = 20251222 }
    @{ Row = 45; E = 10; F = 20251225 }
    @{ Row = 46; E = 4; F = 20251222 }
    @{ Row = 47; E = 3; F = 20251218 }
    @{ Row = 48; E = 4; F = 20251222 }
    @{ Row = 49; E = 5; F = 20251223 }
    @{ Row = 50; E = 8; F = 20251223 }
    @{ Row = 51; E = 8; F = 20251223 }
    @{ Row = 52; E = 8; F = 20251223 }
    @{ Row = 53; E = 8; F = 20251223 }
    @{ Row = 54; E = 8; F = 20251223 }
    @{ Row = 55; E = 8; F = 20251223 }
    @{ Row = 56; E = 8; F = 20251223 }
    @{ Row = 57; E = 8; F = 20251223 }
    @{ Row = 58; E = 2; F = 20251217 }
    @{ Row = 59; E = 2; F = 20251217 }
    @{ Row = 60; E = 2; F = 20251217 }
    @{ Row = 61; E = 5; F = 20251223 }
    @{ Row = 62; E = 2; F = 20251217 }
    @{ Row = 63; E = 2; F = 20251217 }
    @{ Row = 64; E = 2; F = 20251217 }
    @{ Row = 65; E = 3; F = 20251218 }
    @{ Row = 66; E = 3; F = 20251218 }
    @{ Row = 67; E = 3; F = 20251218 }
    @{ Row = 68; E = 3; F = 20251218 }
    @{ Row = 69; E = 3; F = 20251218 }
    @{ Row = 70; E = 4; F = 20251219 }
    @{ Row = 71; E = 4; F = 20251219 }
    @{ Row = 72; E = 4; F = 20251219 }
    @{ Row = 73; E = 4; F = 20251219 }
    @{ Row = 74; E = 4; F = 20251219 }
    @{ Row = 75; E = 4; F = 20251219 }
    @{ Row = 76; E = 4; F = 20251219 }
    @{ Row = 77; E = 7; F = 20251222 }
    @{ Row = 78; E = 7; F = 20251222 }
    @{ Row = 79; E = 7; F = 20251222 }
    @{ Row = 80; E = 7; F = 20251222 }
    @{ Row = 81; E = 7; F = 20251222 }
    @{ Row = 82; E = 7; F = 20251222 }
    @{ Row = 83; E = 7; F = 20251222 }
    @{ Row = 84; E = 7; F = 20251222 }
    @{ Row = 85; E = 7; F = 20251222 }
    @{ Row = 86; E = 7; F = 20251222 }
    @{ Row = 87; E = 4; F = 20251222 }
    @{ Row = 88; E = 4; F = 20251222 }
    @{ Row = 89; E = 4; F = 20251222 }
    @{ Row = 90; E = 4; F = 20251222 }
    @{ Row = 91; E = 10; F = 20251225 }
    @{ Row = 92; E = 4; F = 20251222 }
    @{ Row = 93; E = 7; F = 20251222 }
    @{ Row = 94; E = 7; F = 20251225 }
    @{ Row = 95; E = 6; F = 20251221 }
    @{ Row = 96; E = 4; F = 20251219 }
    @{ Row = 97; E = 4; F = 20251219 }
    @{ Row = 98; E = 4; F = 20251219 }
    @{ Row = 99; E = 4; F = 20251219 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 6).Value = $u.F
}

Write-Host "Updated $($updates.Count) rows"